$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells that hold numeric-looking text (Price / Volume columns) must be
# forced to Text so values like "7.880" or "0.87%" are preserved exactly
# as literal strings instead of being parsed into numbers/percentages.
function Set-TextCell($ws, $addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextCell $ws "D2" "296.78"
Set-TextCell $ws "E2" "0.87%"
Set-TextCell $ws "D3" "31.68"
Set-TextCell $ws "E3" "2.13%"
Set-TextCell $ws "D4" "4.962"
Set-TextCell $ws "E4" "0.67%"
Set-TextCell $ws "D5" "0.07638"
Set-TextCell $ws "E5" "4.00%"
Set-TextCell $ws "D6" "2.252"
Set-TextCell $ws "E6" "-2.21%"
Set-TextCell $ws "D7" "7.880"
Set-TextCell $ws "E7" "1.95%"
$ws.Range("B8").Value = "GateToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextCell $ws "D8" "3.792"
Set-TextCell $ws "E8" "1.37%"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextCell $ws "D9" "0.9251"
Set-TextCell $ws "E9" "2.26%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextCell $ws "D10" "0.09677"
Set-TextCell $ws "E10" "21.43%"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextCell $ws "D11" "0.1737"
Set-TextCell $ws "E11" "3.13%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextCell $ws "D12" "0.08390"
Set-TextCell $ws "E12" "3.18%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextCell $ws "D13" "0.03270"
Set-TextCell $ws "E13" "5.55%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextCell $ws "D14" "0.09838"
Set-TextCell $ws "E14" "-2.53%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextCell $ws "D15" "0.001472"
Set-TextCell $ws "E15" "-2.85%"
$ws.Range("B16").Value = "CoinExToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextCell $ws "D16" "0.04524"
Set-TextCell $ws "E16" "-0.07%"
$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextCell $ws "D17" "0.005761"
Set-TextCell $ws "E17" "-0.92%"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextCell $ws "D18" "3.501"
Set-TextCell $ws "E18" "0.59%"
$ws.Range("B19").Value = "BTSEToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextCell $ws "D19" "2.196"
Set-TextCell $ws "E19" "5.89%"
$ws.Range("B20").Value = "BitpandaEcosystemToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
Set-TextCell $ws "D20" "0.3354"
Set-TextCell $ws "E20" "0.70%"
$ws.Range("B21").Value = "ProBitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
Set-TextCell $ws "D21" "0.1322"
Set-TextCell $ws "E21" "1.50%"
$ws.Range("B22").Value = "MCDex"
$ws.Range("C22").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextCell $ws "D22" "4.072"
Set-TextCell $ws "E22" "2.59%"
$ws.Range("B23").Value = "ZBToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
Set-TextCell $ws "D23" "0.2280"
Set-TextCell $ws "E23" "8.76%"
Set-TextCell $ws "D24" "0.001214"
Set-TextCell $ws "E24" "0.26%"
Set-TextCell $ws "D25" "0.004343"
Set-TextCell $ws "E25" "-6.74%"
Set-TextCell $ws "D26" "0.0001292"
Set-TextCell $ws "E26" "-0.64%"
Set-TextCell $ws "D27" "0.0003373"
Set-TextCell $ws "E27" "-0.49%"
Set-TextCell $ws "D39" "0.01675"
Set-TextCell $ws "E39" "3.84%"
Set-TextCell $ws "E40" "4.04%"
Set-TextCell $ws "D41" "0.007482"
Set-TextCell $ws "E41" "1.87%"
Set-TextCell $ws "D42" "0.009738"
Set-TextCell $ws "E42" "13.74%"
Set-TextCell $ws "D43" "0.1382"
Set-TextCell $ws "E43" "3.96%"
Set-TextCell $ws "D44" "0.002101"
Set-TextCell $ws "E44" "5.06%"
Set-TextCell $ws "D45" "0.009432"
Set-TextCell $ws "E45" "-0.71%"
Set-TextCell $ws "D46" "0.00006070"
Set-TextCell $ws "E46" "2.56%"
Set-TextCell $ws "E47" "-0.48%"
Set-TextCell $ws "D48" "2.551"
Set-TextCell $ws "E48" "13.83%"
Set-TextCell $ws "D49" "0.001987"
Set-TextCell $ws "E49" "-31.33%"
Set-TextCell $ws "E50" "-0.48%"
Set-TextCell $ws "E51" "-0.48%"
